$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the grade values in row 14
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 5
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 5

# Move the active selection to H14 (updates sheetView pane/selection)
$ws.Range("H14").Select()
